$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(3).TextFrame.TextRange.Text = "Trigger (Sentinel Incident, Sentinel Alert, Sentinel Entity, MDC Recommendation, MDC Alert, MDC Compliance)"

$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Enrich – Get a geolocation, check against Virus Total or URL/IP reputation, Look up other fields of data on entity, is that user in a sales role who travels ? Is OOF Turned on ?"

# Refresh the cached "datetimeFigureOut" placeholder text (slide master,
# every custom layout, and the notes master) to the current save date,
# mirroring PowerPoint's automatic field recache on re-save.
$newDate = "12/12/2023"

$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$layoutDateIdx = @{1=3; 2=3; 3=3; 4=4; 5=6; 6=2; 7=1; 8=4; 9=4; 10=3; 11=3}
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
  $layout = $layouts.Item($li)
  $dateIdx = $layoutDateIdx[$li]
  $layout.Shapes.Item($dateIdx).TextFrame.TextRange.Text = $newDate
}

$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate
